# Applies the textual updates described in the commit diff:
# - Updates the date heading from 2024-12-30 Monday to 2024-12-31 Tuesday
# - Updates each arithmetic expression cell in the table to its new value
#
# Each "old" value is unique within the document, so a simple
# Find/Replace (scoped to the whole document content) for each pair is
# sufficient and unambiguous.

$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-12-30 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-12-31 Tuesday", 2) | Out-Null
$d.Content.Find.Execute("1+79=", $true, $false, $false, $false, $false, $true, 1, $false, "91-70=", 2) | Out-Null
$d.Content.Find.Execute("90-24=", $true, $false, $false, $false, $false, $true, 1, $false, "2+47=", 2) | Out-Null
$d.Content.Find.Execute("14+46=", $true, $false, $false, $false, $false, $true, 1, $false, "59+3=", 2) | Out-Null
$d.Content.Find.Execute("61-49=", $true, $false, $false, $false, $false, $true, 1, $false, "19+49=", 2) | Out-Null
$d.Content.Find.Execute("9-6=", $true, $false, $false, $false, $false, $true, 1, $false, "80-2=", 2) | Out-Null
$d.Content.Find.Execute("17+52=", $true, $false, $false, $false, $false, $true, 1, $false, "67+27=", 2) | Out-Null
$d.Content.Find.Execute("90-89=", $true, $false, $false, $false, $false, $true, 1, $false, "94-71=", 2) | Out-Null
$d.Content.Find.Execute("60-36=", $true, $false, $false, $false, $false, $true, 1, $false, "11+43=", 2) | Out-Null
$d.Content.Find.Execute("82-3=", $true, $false, $false, $false, $false, $true, 1, $false, "98-42=", 2) | Out-Null
$d.Content.Find.Execute("78-58=", $true, $false, $false, $false, $false, $true, 1, $false, "6+3=", 2) | Out-Null
$d.Content.Find.Execute("90+9=", $true, $false, $false, $false, $false, $true, 1, $false, "92-26=", 2) | Out-Null
$d.Content.Find.Execute("34+18=", $true, $false, $false, $false, $false, $true, 1, $false, "26+43=", 2) | Out-Null
$d.Content.Find.Execute("97-72=", $true, $false, $false, $false, $false, $true, 1, $false, "37+24=", 2) | Out-Null
$d.Content.Find.Execute("47-38=", $true, $false, $false, $false, $false, $true, 1, $false, "44-30=", 2) | Out-Null
$d.Content.Find.Execute("70-10=", $true, $false, $false, $false, $false, $true, 1, $false, "27+5=", 2) | Out-Null
$d.Content.Find.Execute("96-46=", $true, $false, $false, $false, $false, $true, 1, $false, "38-28=", 2) | Out-Null
$d.Content.Find.Execute("24+41=", $true, $false, $false, $false, $false, $true, 1, $false, "55+8=", 2) | Out-Null
$d.Content.Find.Execute("73-44=", $true, $false, $false, $false, $false, $true, 1, $false, "58-50=", 2) | Out-Null
$d.Content.Find.Execute("88-23=", $true, $false, $false, $false, $false, $true, 1, $false, "93-50=", 2) | Out-Null
$d.Content.Find.Execute("52+37=", $true, $false, $false, $false, $false, $true, 1, $false, "39+8=", 2) | Out-Null
$d.Content.Find.Execute("3+71=", $true, $false, $false, $false, $false, $true, 1, $false, "24+13=", 2) | Out-Null
$d.Content.Find.Execute("96-5=", $true, $false, $false, $false, $false, $true, 1, $false, "45+45=", 2) | Out-Null
$d.Content.Find.Execute("45+41=", $true, $false, $false, $false, $false, $true, 1, $false, "8+0=", 2) | Out-Null
$d.Content.Find.Execute("1+55=", $true, $false, $false, $false, $false, $true, 1, $false, "6+45=", 2) | Out-Null
$d.Content.Find.Execute("42-20=", $true, $false, $false, $false, $false, $true, 1, $false, "70+29=", 2) | Out-Null
$d.Content.Find.Execute("57-18=", $true, $false, $false, $false, $false, $true, 1, $false, "49-40=", 2) | Out-Null
$d.Content.Find.Execute("87-11=", $true, $false, $false, $false, $false, $true, 1, $false, "59-33=", 2) | Out-Null
$d.Content.Find.Execute("79+9=", $true, $false, $false, $false, $false, $true, 1, $false, "11-4=", 2) | Out-Null
$d.Content.Find.Execute("44+10=", $true, $false, $false, $false, $false, $true, 1, $false, "5+86=", 2) | Out-Null
$d.Content.Find.Execute("46-34=", $true, $false, $false, $false, $false, $true, 1, $false, "30+5=", 2) | Out-Null
$d.Content.Find.Execute("64-61=", $true, $false, $false, $false, $false, $true, 1, $false, "38-24=", 2) | Out-Null
$d.Content.Find.Execute("93-33=", $true, $false, $false, $false, $false, $true, 1, $false, "59+27=", 2) | Out-Null
$d.Content.Find.Execute("42-14=", $true, $false, $false, $false, $false, $true, 1, $false, "26-2=", 2) | Out-Null
$d.Content.Find.Execute("55+31=", $true, $false, $false, $false, $false, $true, 1, $false, "49+37=", 2) | Out-Null
$d.Content.Find.Execute("83-82=", $true, $false, $false, $false, $false, $true, 1, $false, "82-6=", 2) | Out-Null
$d.Content.Find.Execute("82-65=", $true, $false, $false, $false, $false, $true, 1, $false, "65+27=", 2) | Out-Null
$d.Content.Find.Execute("76+10=", $true, $false, $false, $false, $false, $true, 1, $false, "5+11=", 2) | Out-Null
$d.Content.Find.Execute("93-0=", $true, $false, $false, $false, $false, $true, 1, $false, "58+29=", 2) | Out-Null
$d.Content.Find.Execute("85-39=", $true, $false, $false, $false, $false, $true, 1, $false, "98-42=", 2) | Out-Null
$d.Content.Find.Execute("77-16=", $true, $false, $false, $false, $false, $true, 1, $false, "24+69=", 2) | Out-Null
$d.Content.Find.Execute("16+58=", $true, $false, $false, $false, $false, $true, 1, $false, "90-67=", 2) | Out-Null
$d.Content.Find.Execute("52+3=", $true, $false, $false, $false, $false, $true, 1, $false, "6-2=", 2) | Out-Null
$d.Content.Find.Execute("54-1=", $true, $false, $false, $false, $false, $true, 1, $false, "82+5=", 2) | Out-Null
$d.Content.Find.Execute("81-31=", $true, $false, $false, $false, $false, $true, 1, $false, "49-31=", 2) | Out-Null
$d.Content.Find.Execute("81-19=", $true, $false, $false, $false, $false, $true, 1, $false, "75-54=", 2) | Out-Null
$d.Content.Find.Execute("51+44=", $true, $false, $false, $false, $false, $true, 1, $false, "63-40=", 2) | Out-Null
$d.Content.Find.Execute("55+44=", $true, $false, $false, $false, $false, $true, 1, $false, "86-34=", 2) | Out-Null
$d.Content.Find.Execute("36+6=", $true, $false, $false, $false, $false, $true, 1, $false, "57+37=", 2) | Out-Null
$d.Content.Find.Execute("79-7=", $true, $false, $false, $false, $false, $true, 1, $false, "56-35=", 2) | Out-Null
$d.Content.Find.Execute("29-26=", $true, $false, $false, $false, $false, $true, 1, $false, "49-40=", 2) | Out-Null
$d.Content.Find.Execute("41-2=", $true, $false, $false, $false, $false, $true, 1, $false, "59-1=", 2) | Out-Null
$d.Content.Find.Execute("55+35=", $true, $false, $false, $false, $false, $true, 1, $false, "29+26=", 2) | Out-Null
$d.Content.Find.Execute("51-11=", $true, $false, $false, $false, $false, $true, 1, $false, "56-2=", 2) | Out-Null
$d.Content.Find.Execute("1+16=", $true, $false, $false, $false, $false, $true, 1, $false, "81-14=", 2) | Out-Null
$d.Content.Find.Execute("69-10=", $true, $false, $false, $false, $false, $true, 1, $false, "94-55=", 2) | Out-Null
$d.Content.Find.Execute("97-85=", $true, $false, $false, $false, $false, $true, 1, $false, "12+37=", 2) | Out-Null
$d.Content.Find.Execute("42+5=", $true, $false, $false, $false, $false, $true, 1, $false, "31-18=", 2) | Out-Null
$d.Content.Find.Execute("37+21=", $true, $false, $false, $false, $false, $true, 1, $false, "46-37=", 2) | Out-Null
$d.Content.Find.Execute("18+80=", $true, $false, $false, $false, $false, $true, 1, $false, "37-32=", 2) | Out-Null
$d.Content.Find.Execute("62-22=", $true, $false, $false, $false, $false, $true, 1, $false, "72-69=", 2) | Out-Null
$d.Content.Find.Execute("36+14=", $true, $false, $false, $false, $false, $true, 1, $false, "26-13=", 2) | Out-Null
$d.Content.Find.Execute("41+47=", $true, $false, $false, $false, $false, $true, 1, $false, "56+43=", 2) | Out-Null
$d.Content.Find.Execute("84+0=", $true, $false, $false, $false, $false, $true, 1, $false, "28+27=", 2) | Out-Null
$d.Content.Find.Execute("17+34=", $true, $false, $false, $false, $false, $true, 1, $false, "16+29=", 2) | Out-Null
$d.Content.Find.Execute("12+87=", $true, $false, $false, $false, $false, $true, 1, $false, "85-23=", 2) | Out-Null
$d.Content.Find.Execute("74-30=", $true, $false, $false, $false, $false, $true, 1, $false, "23+19=", 2) | Out-Null
$d.Content.Find.Execute("66-29=", $true, $false, $false, $false, $false, $true, 1, $false, "41+22=", 2) | Out-Null
$d.Content.Find.Execute("39+15=", $true, $false, $false, $false, $false, $true, 1, $false, "72-49=", 2) | Out-Null
$d.Content.Find.Execute("79+13=", $true, $false, $false, $false, $false, $true, 1, $false, "44-5=", 2) | Out-Null
$d.Content.Find.Execute("86-18=", $true, $false, $false, $false, $false, $true, 1, $false, "37+13=", 2) | Out-Null
$d.Content.Find.Execute("33-5=", $true, $false, $false, $false, $false, $true, 1, $false, "29+69=", 2) | Out-Null
$d.Content.Find.Execute("24-2=", $true, $false, $false, $false, $false, $true, 1, $false, "76+8=", 2) | Out-Null
$d.Content.Find.Execute("74+23=", $true, $false, $false, $false, $false, $true, 1, $false, "9+77=", 2) | Out-Null
$d.Content.Find.Execute("44+38=", $true, $false, $false, $false, $false, $true, 1, $false, "2+21=", 2) | Out-Null
$d.Content.Find.Execute("40-10=", $true, $false, $false, $false, $false, $true, 1, $false, "72-42=", 2) | Out-Null
$d.Content.Find.Execute("92-5=", $true, $false, $false, $false, $false, $true, 1, $false, "72+4=", 2) | Out-Null
$d.Content.Find.Execute("19+31=", $true, $false, $false, $false, $false, $true, 1, $false, "44-13=", 2) | Out-Null
$d.Content.Find.Execute("34-7=", $true, $false, $false, $false, $false, $true, 1, $false, "50-15=", 2) | Out-Null
$d.Content.Find.Execute("64-36=", $true, $false, $false, $false, $false, $true, 1, $false, "16+71=", 2) | Out-Null
$d.Content.Find.Execute("62+34=", $true, $false, $false, $false, $false, $true, 1, $false, "91-59=", 2) | Out-Null
$d.Content.Find.Execute("31+68=", $true, $false, $false, $false, $false, $true, 1, $false, "8+74=", 2) | Out-Null
$d.Content.Find.Execute("50-37=", $true, $false, $false, $false, $false, $true, 1, $false, "49+11=", 2) | Out-Null
$d.Content.Find.Execute("96-25=", $true, $false, $false, $false, $false, $true, 1, $false, "42-36=", 2) | Out-Null
$d.Content.Find.Execute("80+9=", $true, $false, $false, $false, $false, $true, 1, $false, "90-58=", 2) | Out-Null
$d.Content.Find.Execute("97-25=", $true, $false, $false, $false, $false, $true, 1, $false, "94-67=", 2) | Out-Null
$d.Content.Find.Execute("86-68=", $true, $false, $false, $false, $false, $true, 1, $false, "76-32=", 2) | Out-Null
$d.Content.Find.Execute("44-21=", $true, $false, $false, $false, $false, $true, 1, $false, "40+51=", 2) | Out-Null
$d.Content.Find.Execute("43-25=", $true, $false, $false, $false, $false, $true, 1, $false, "3+65=", 2) | Out-Null
$d.Content.Find.Execute("65+20=", $true, $false, $false, $false, $false, $true, 1, $false, "8+88=", 2) | Out-Null
$d.Content.Find.Execute("37+19=", $true, $false, $false, $false, $false, $true, 1, $false, "35+62=", 2) | Out-Null
$d.Content.Find.Execute("44-2=", $true, $false, $false, $false, $false, $true, 1, $false, "97-17=", 2) | Out-Null
$d.Content.Find.Execute("18-0=", $true, $false, $false, $false, $false, $true, 1, $false, "11+79=", 2) | Out-Null
$d.Content.Find.Execute("78-57=", $true, $false, $false, $false, $false, $true, 1, $false, "73-60=", 2) | Out-Null
$d.Content.Find.Execute("96-26=", $true, $false, $false, $false, $false, $true, 1, $false, "59-4=", 2) | Out-Null
$d.Content.Find.Execute("51+17=", $true, $false, $false, $false, $false, $true, 1, $false, "84+2=", 2) | Out-Null
$d.Content.Find.Execute("22+62=", $true, $false, $false, $false, $false, $true, 1, $false, "8+38=", 2) | Out-Null
$d.Content.Find.Execute("48-45=", $true, $false, $false, $false, $false, $true, 1, $false, "90-55=", 2) | Out-Null
$d.Content.Find.Execute("62-54=", $true, $false, $false, $false, $false, $true, 1, $false, "34-22=", 2) | Out-Null
$d.Content.Find.Execute("64-9=", $true, $false, $false, $false, $false, $true, 1, $false, "69-31=", 2) | Out-Null
$d.Content.Find.Execute("36-5=", $true, $false, $false, $false, $false, $true, 1, $false, "59-42=", 2) | Out-Null

Write-Output "Applied 101 replacements"
